$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: Morning Glass of Ether
$ws.Range("H15").Value = 92.56999999999999
$ws.Range("I15").Value = 92.56999999999999
$ws.Range("K15").Value = 277.71
$ws.Range("M15").Value = -108.71

# Row 42: Eye of the Beholder
$ws.Range("H42").Value = 215.08333
$ws.Range("J42").Value = 170
$ws.Range("L42").Value = 510
$ws.Range("N42").Value = -970

# Row 111: An Eye for Healing
$ws.Range("H111").Value = 1005.7692
$ws.Range("I111").Value = 717.5
$ws.Range("J111").Value = 1966.6666
$ws.Range("K111").Value = 2152.5
$ws.Range("L111").Value = 5899.9998
$ws.Range("M111").Value = 914.5
$ws.Range("N111").Value = -12033.9998

# Row 139: Something Salty and Ceremonial
$ws.Range("H139").Value = 70337.5
$ws.Range("J139").Value = 70337.5
$ws.Range("L139").Value = 70337.5
$ws.Range("N139").Value = -80617.5

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Range("H32").Value = 182935.58
$ws.Range("I32").Value = 7013.5576
$ws.Range("J32").Value = 853638.25
$ws.Range("K32").Value = 7013.5576
$ws.Range("L32").Value = 853638.25
$ws.Range("M32").Value = -6726.5576
$ws.Range("N32").Value = -854212.25

# Row 37: Get Shirty
$ws.Range("H37").Value = 7143.2856
$ws.Range("J37").Value = 8000.75
$ws.Range("L37").Value = 8000.75
$ws.Range("N37").Value = -8546.75

# Row 61: Dealing with the Tough Stuff
$ws.Range("H61").Value = 3295.2856
$ws.Range("J61").Value = 3022.1
$ws.Range("L61").Value = 3022.1
$ws.Range("N61").Value = -3446.1

# Row 74: As the Bolt Flies
$ws.Range("H74").Value = 970.36365
$ws.Range("I74").Value = 720
$ws.Range("J74").Value = 1638
$ws.Range("K74").Value = 720
$ws.Range("L74").Value = 1638
$ws.Range("M74").Value = 154
$ws.Range("N74").Value = -3386

# Row 77: Heavy Metal Banned (L)
$ws.Range("H77").Value = 970.36365
$ws.Range("I77").Value = 720
$ws.Range("J77").Value = 1638
$ws.Range("K77").Value = 3600
$ws.Range("L77").Value = 8190
$ws.Range("M77").Value = 768
$ws.Range("N77").Value = -16926

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 4813.6875
$ws.Range("I122").Value = 5014.0435
$ws.Range("K122").Value = 15042.1305
$ws.Range("M122").Value = -12592.1305

# Row 132: Don't Bore Me, Ore Me
$ws.Range("H132").Value = 2058.2778
$ws.Range("I132").Value = 1837.6154
$ws.Range("J132").Value = 2632
$ws.Range("K132").Value = 5512.8462
$ws.Range("L132").Value = 7896
$ws.Range("M132").Value = -2982.8462
$ws.Range("N132").Value = -12956

# Row 136: Metal with Mettle
$ws.Range("H136").Value = 3295.2856
$ws.Range("J136").Value = 3022.1
$ws.Range("L136").Value = 9066.299999999999
$ws.Range("N136").Value = -14166.3

$ws = $wb.Worksheets.Item("BSM")
# Row 20: Smelt and Dealt
$ws.Range("H20").Value = 3714.1667
$ws.Range("I20").Value = 3242.4
$ws.Range("J20").Value = 4051.1428
$ws.Range("K20").Value = 3242.4
$ws.Range("L20").Value = 4051.1428
$ws.Range("M20").Value = -2995.4
$ws.Range("N20").Value = -4545.1428

# Row 105: Ingot to Wing It
$ws.Range("H105").Value = 2000
$ws.Range("I105").Value = 2000
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2000
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -253

# Row 134: Ruthenium Supremium
$ws.Range("H134").Value = 6400.154
$ws.Range("I134").Value = 900.1667
$ws.Range("K134").Value = 2700.5001
$ws.Range("M134").Value = -165.5001000000002

# Row 140: Ceremonial Teeth
$ws.Range("H140").Value = 89200
$ws.Range("J140").Value = 89200
$ws.Range("L140").Value = 89200
$ws.Range("N140").Value = -99560

$ws = $wb.Worksheets.Item("CRP")
# Row 60: Bowing to Greater Power
$ws.Range("H60").Value = 7279.2
$ws.Range("J60").Value = 8101
$ws.Range("L60").Value = 8101
$ws.Range("N60").Value = -9123

# Row 68: Do You Even String Bow
$ws.Range("H68").Value = 14537
$ws.Range("J68").Value = 14537
$ws.Range("L68").Value = 14537
$ws.Range("N68").Value = -16035

# Row 71: Win One Bow, Get Three Free (L)
$ws.Range("H71").Value = 14537
$ws.Range("J71").Value = 14537
$ws.Range("L71").Value = 43611
$ws.Range("N71").Value = -51099

# Row 74: License to Heal
$ws.Range("H74").Value = 18828
$ws.Range("J74").Value = 18828
$ws.Range("L74").Value = 18828
$ws.Range("N74").Value = -20576

# Row 77: Purified Polyrhythm (L)
$ws.Range("H77").Value = 18828
$ws.Range("J77").Value = 18828
$ws.Range("L77").Value = 56484
$ws.Range("N77").Value = -65220

# Row 105: Zelkova, My Love
$ws.Range("H105").Value = 1041.25
$ws.Range("I105").Value = 704.61536
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 704.61536
$ws.Range("L105").Value = 2500
$ws.Range("M105").Value = 1042.38464
$ws.Range("N105").Value = -5994

# Row 122: Timber of Tenkonto
$ws.Range("H122").Value = 1649.7693
$ws.Range("I122").Value = 1476.2
$ws.Range("J122").Value = 2228.3333
$ws.Range("K122").Value = 4428.6
$ws.Range("L122").Value = 6684.999899999999
$ws.Range("M122").Value = -1978.6
$ws.Range("N122").Value = -11584.9999

# Row 138: Bow Out
$ws.Range("H138").Value = 49275
$ws.Range("J138").Value = 49275
$ws.Range("L138").Value = 49275
$ws.Range("N138").Value = -59555

$ws = $wb.Worksheets.Item("CUL")
# Row 121: A Cookie for Your Troubles
$ws.Range("H121").Value = 74665.78
$ws.Range("I121").Value = 315.4
$ws.Range("J121").Value = 91563.59
$ws.Range("K121").Value = 946.1999999999999
$ws.Range("L121").Value = 274690.77
$ws.Range("M121").Value = 363.8000000000001
$ws.Range("N121").Value = -277310.77

# Row 122: Salt of the North
$ws.Range("H122").Value = 584.02325
$ws.Range("I122").Value = 352.9375
$ws.Range("J122").Value = 1256.2727
$ws.Range("K122").Value = 3176.4375
$ws.Range("L122").Value = 11306.4543
$ws.Range("M122").Value = -726.4375
$ws.Range("N122").Value = -16206.4543

# Row 131: The Mountain Steeped
$ws.Range("H131").Value = 6098431.5
$ws.Range("I131").Value = 962.6
$ws.Range("J131").Value = 8065357
$ws.Range("K131").Value = 2887.8
$ws.Range("L131").Value = 24196071
$ws.Range("M131").Value = 2152.2
$ws.Range("N131").Value = -24206151

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 557002
$ws.Range("I80").Value = 557002
$ws.Range("K80").Value = 557002
$ws.Range("M80").Value = -556004

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 557002
$ws.Range("I83").Value = 557002
$ws.Range("K83").Value = 2785010
$ws.Range("M83").Value = -2780018

# Row 122: Awarding Academic Excellence
$ws.Range("H122").Value = 2507.4167
$ws.Range("I122").Value = 2235.9565
$ws.Range("J122").Value = 2987.6924
$ws.Range("K122").Value = 6707.869499999999
$ws.Range("L122").Value = 8963.0772
$ws.Range("M122").Value = -4257.869499999999
$ws.Range("N122").Value = -13863.0772

# Row 123: Workplace Workout
$ws.Range("H123").Value = 12392.759
$ws.Range("J123").Value = 12392.759
$ws.Range("L123").Value = 12392.759
$ws.Range("N123").Value = -17292.759

# Row 140: The Right Rod
$ws.Range("H140").Value = 89749.5
$ws.Range("J140").Value = 89749.5
$ws.Range("L140").Value = 89749.5
$ws.Range("N140").Value = -100109.5

$ws = $wb.Worksheets.Item("LTW")
# Row 36: Campaign in the Membrane
$ws.Range("H36").Value = 37500
$ws.Range("J36").Value = 37500
$ws.Range("L36").Value = 37500
$ws.Range("N36").Value = -38624

# Row 40: Best Served Toad
$ws.Range("H40").Value = 2966.72
$ws.Range("I40").Value = 2791.889
$ws.Range("J40").Value = 3416.2856
$ws.Range("K40").Value = 2791.889
$ws.Range("L40").Value = 3416.2856
$ws.Range("M40").Value = -2655.889
$ws.Range("N40").Value = -3688.2856

# Row 61: Spelling Me Softly
$ws.Range("H61").Value = 2241.2778
$ws.Range("I61").Value = 1734.3846
$ws.Range("J61").Value = 3559.2
$ws.Range("K61").Value = 1734.3846
$ws.Range("L61").Value = 3559.2
$ws.Range("M61").Value = -1532.3846
$ws.Range("N61").Value = -3963.2

# Row 93: Hide to Go Seek
$ws.Range("H93").Value = 2556.1304
$ws.Range("I93").Value = 1846
$ws.Range("J93").Value = 3102.3845
$ws.Range("K93").Value = 1846
$ws.Range("L93").Value = 3102.3845
$ws.Range("M93").Value = -598
$ws.Range("N93").Value = -5598.3845

# Row 113: Peace in Rest
$ws.Range("H113").Value = 2241.2778
$ws.Range("I113").Value = 1734.3846
$ws.Range("J113").Value = 3559.2
$ws.Range("K113").Value = 1734.3846
$ws.Range("L113").Value = 3559.2
$ws.Range("M113").Value = 435.6153999999999
$ws.Range("N113").Value = -7899.2

# Row 122: Hell on Leather
$ws.Range("H122").Value = 2553.0454
$ws.Range("I122").Value = 2380
$ws.Range("J122").Value = 4283.5
$ws.Range("K122").Value = 7140
$ws.Range("L122").Value = 12850.5
$ws.Range("M122").Value = -4690
$ws.Range("N122").Value = -17750.5

# Row 138: Freezing Toes
$ws.Range("H138").Value = 61306.332
$ws.Range("J138").Value = 61306.332
$ws.Range("L138").Value = 61306.332
$ws.Range("N138").Value = -71586.33199999999

# Row 139: Giving Gatherers Their Gear
$ws.Range("H139").Value = 56280
$ws.Range("J139").Value = 69100
$ws.Range("L139").Value = 69100
$ws.Range("N139").Value = -79380
